$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first 16 data rows (rows 2 through 17), shifting the
# remaining data up. This removes the earliest 16 revision/date
# observations from the series.
$ws.Range("A2:B17").EntireRow.Delete()
